$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STORAGE")

for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Range("B$row")
    if ($cell.Value2 -eq "ELECTRIC_HEAVY_GOODS_VEHICLE") {
        $cell.Value = "ELECTRIC_VEHICLE"
    }
}

$wb.Save()
